$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 85
$ws.Range("I2").Value = 311
$ws.Range("J2").Value = 1255
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 337
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = 206
$ws.Range("P2").Value = 4
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 140
$ws.Range("T2").Value = 219
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 1800
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1871
$ws.Range("Z2").Value = 24
$ws.Range("AA2").Value = 8
